$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in AD1:AF1, matching the
# existing header formatting (bold, centered, bordered) by copying the
# format from an existing header cell first, then overwriting the values.
$ws.Range("A1").Copy($ws.Range("AD1:AF1"))
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row 2-44.
$ws.Range("AD2:AD44").Value = 85
$ws.Range("AE2:AE44").Value = 77
$ws.Range("AF2:AF44").Value = 0
